$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos -> short bio text replaces long PT objectives paragraph ---
$ws.Range("B10").Value = "5817330 - Larissa de Freitas"
$ws.Range("C10").Value = "5817330 - Larissa de Freitas"

# --- Row 13: was "5817330 - Larissa de Freitas" -> becomes "Programa resumido:" / "01/01/2022" ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2022"
$ws.Range("C13").Value = "01/01/2022"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: was "1506103 - Pedro Carlos de Oliveira" -> becomes "Short syllabus:" / short EN syllabus ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Chemical Kinetics, Thermodynamics and Equilibrium, Chemical Equilibrium, Eletrochemistry."
$ws.Range("C14").Value = "Chemical Kinetics, Thermodynamics and Equilibrium, Chemical Equilibrium, Eletrochemistry."
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: was "Programa resumido:" / short PT syllabus -> becomes "Programa:" / bio text ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "5817330 - Larissa de Freitas"
$ws.Range("C15").Value = "5817330 - Larissa de Freitas"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: was "Short syllabus:" / short EN syllabus -> becomes "Syllabus:" / long EN syllabus ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Chemical Kinetics: Effect of concentration of reactants on the rate of chemical reactions. Equations reaction speed differentials. Rate laws for 1st and 2nd order reactions. Activation energy and catalysts. Arrhenius equation. Thermodynamics and equilibrium: First Law of Thermodynamics. Reaction heat and internal energy. Variation of enthalpy and enthalpy. Entropy. Second Law of Thermodynamics. Entropy and disorder. Third Law of Thermodynamics. Entropy variation of a reaction. The Standard Free Energy and spontaneity criterion. Relationship between standard free energy change and equilibrium constant. Chemical Equilibrium: Nature's chemical balance. Quotient reaction and equilibrium constant. Effect of reagent concentration and temperature on the equilibrium (Le Chatelier's principle). Equilibrium Acid-base. Equilibrium in solution: solubility and complex ions. Electrochemistry: Semi-reactions. Standard electrode potential. Potential galvanic cells. Relationship between standard free energy change and cell potential. Gibbs free energy and Nernst equation. Electrolysis and Faraday law."
$ws.Range("C16").Value = "Chemical Kinetics: Effect of concentration of reactants on the rate of chemical reactions. Equations reaction speed differentials. Rate laws for 1st and 2nd order reactions. Activation energy and catalysts. Arrhenius equation. Thermodynamics and equilibrium: First Law of Thermodynamics. Reaction heat and internal energy. Variation of enthalpy and enthalpy. Entropy. Second Law of Thermodynamics. Entropy and disorder. Third Law of Thermodynamics. Entropy variation of a reaction. The Standard Free Energy and spontaneity criterion. Relationship between standard free energy change and equilibrium constant. Chemical Equilibrium: Nature's chemical balance. Quotient reaction and equilibrium constant. Effect of reagent concentration and temperature on the equilibrium (Le Chatelier's principle). Equilibrium Acid-base. Equilibrium in solution: solubility and complex ions. Electrochemistry: Semi-reactions. Standard electrode potential. Potential galvanic cells. Relationship between standard free energy change and cell potential. Gibbs free energy and Nernst equation. Electrolysis and Faraday law."
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: was "Programa:" / long PT program -> becomes "Avaliação:" (A only, B/C cleared) ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Rows.Item(17).RowHeight = 15

# --- Row 18: was "Syllabus:" / long EN syllabus -> becomes "Método:" / bio text ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1506103 - Pedro Carlos de Oliveira"
$ws.Range("C18").Value = "1506103 - Pedro Carlos de Oliveira"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: was "Avaliação:" (A only) -> becomes "Critério:" / exam description ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Serão realizadas duas provas escritas"
$ws.Range("C19").Value = "Serão realizadas duas provas escritas"
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: was "Método:" / exam description -> becomes "Norma de recuperação:" / NF formula ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "NF = (P1 + P2*2)/3"
$ws.Range("C20").Value = "NF = (P1 + P2*2)/3"
$ws.Rows.Item(20).RowHeight = 15

# --- Row 21: was "Critério:" / NF formula -> becomes "Bibliografia:" / recovery norm text ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Será realizada uma avaliação (P3) englobando toda a ementa. A média final será obtida conforme equação: MF= (NF+P3)/2."
$ws.Range("C21").Value = "Será realizada uma avaliação (P3) englobando toda a ementa. A média final será obtida conforme equação: MF= (NF+P3)/2."
$ws.Rows.Item(21).RowHeight = 120

# --- Rows 22 and 23 are removed entirely (old recovery-norm text row and bibliography row) ---
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
